# "Generate Report for Archive"
# The localization status moved on from "Ready for handoff" to "In Translation"
# for this file's zh-cn / de-de targets. Update every place that status string
# appears (the Overview roll-up sheet and each language-specific detail sheet),
# then re-fit the status column(s) now that the text is shorter.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: status is duplicated per-language in columns E (zh-cn) and F (de-de), row 2
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsOverview.Range("E2:F2").EntireColumn.AutoFit() | Out-Null

# --- zh-cn detail sheet: "Status" column is C, row 2
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Range("C2").EntireColumn.AutoFit() | Out-Null

# --- de-de detail sheet: "Status" column is C, row 2
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Range("C2").EntireColumn.AutoFit() | Out-Null
